$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point precision of the timestamp in A47
$ws.Range("A47").Value = 44360.76849045602

# Add new row 48 data
$ws.Range("A48").Value = 44361.76756595966
$ws.Range("B48").Value = 77087
$ws.Range("C48").Value = 64761
$ws.Range("D48").Value = 3327
$ws.Range("E48").Value = 2100
$ws.Range("F48").Value = 1481
$ws.Range("G48").Value = 20381
$ws.Range("H48").Value = 1409
$ws.Range("I48").Value = 884
$ws.Range("J48").Value = 177
